# Update DailyStats rows with revised AgTests (F) / AgPosit (G) figures
# per the "Updated: pi 07. 05. 2021" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Column F value, Column G value ($null = leave unchanged)
$updates = @(
    @{ Row = 395; F = 750917; G = 1958 },
    @{ Row = 398; F = 298548; G = $null },
    @{ Row = 400; F = $null;  G = 769 },
    @{ Row = 401; F = 273668; G = $null },
    @{ Row = 402; F = 717226; G = 1388 },
    @{ Row = 405; F = 174045; G = $null },
    @{ Row = 408; F = 303976; G = 836 },
    @{ Row = 409; F = 703321; G = 1002 },
    @{ Row = 412; F = 175982; G = $null },
    @{ Row = 414; F = 146317; G = 555 },
    @{ Row = 415; F = 305201; G = 694 },
    @{ Row = 418; F = 200617; G = 698 },
    @{ Row = 422; F = 293800; G = 634 },
    @{ Row = 423; F = 431466; G = 627 },
    @{ Row = 424; F = 253598; G = 484 },
    @{ Row = 425; F = 136181; G = 537 },
    @{ Row = 426; F = 104312; G = 387 },
    @{ Row = 427; F = 88353;  G = 353 }
)

foreach ($u in $updates) {
    if ($null -ne $u.F) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F   # Column F
    }
    if ($null -ne $u.G) {
        $ws.Cells.Item($u.Row, 7).Value = $u.G   # Column G
    }
}
